$wb = $excel.ActiveWorkbook

# --- New sheet: ID_c4955e1 (SanDisk SSD) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)

$ws3.Cells.Item(1,1).Value = "price"
$ws3.Cells.Item(1,2).Value = "date"
$ws3.Cells.Item(1,3).Value = "product"
$header3 = $ws3.Range("A1:C1")
$header3.Font.Bold = $true
$header3.HorizontalAlignment = -4108
$header3.VerticalAlignment = -4160
$header3.Borders.LineStyle = 1

$ws3.Cells.Item(2,1).Value = 109
$ws3.Cells.Item(2,2).NumberFormat = "@"
$ws3.Cells.Item(2,2).Value = "10/02/2025"
$ws3.Cells.Item(2,3).Value = "SanDisk SSD PLUS 1TB Internal SSD - SATA III 6 Gb/s"

$ws3.Name = "ID_c4955e1"

# --- New sheet: ID_2b6fb5b (John Hardy bracelet) ---
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet2)

$ws4.Cells.Item(1,1).Value = "price"
$ws4.Cells.Item(1,2).Value = "date"
$ws4.Cells.Item(1,3).Value = "product"
$header4 = $ws4.Range("A1:C1")
$header4.Font.Bold = $true
$header4.HorizontalAlignment = -4108
$header4.VerticalAlignment = -4160
$header4.Borders.LineStyle = 1

$ws4.Cells.Item(2,1).Value = 695
$ws4.Cells.Item(2,2).NumberFormat = "@"
$ws4.Cells.Item(2,2).Value = "10/02/2025"
$ws4.Cells.Item(2,3).Value = "John Hardy Women's Legends Naga Gold & Silver Dragon Station Chain Bracelet"

$ws4.Name = "ID_2b6fb5b"
